# Auto-generated Word COM-interop script
# Applies the "Pitch de Apresentacao" edit described by the commit diff:
#  1. Rewrites paragraph 1 (ChatGPT project pitch) with new wording and
#     drops its justified alignment (keeps the first-line indent).
#  2. Inserts two new paragraphs (frameworks/libraries + revenue model)
#     right after the "Acreditamos..." paragraph.

$d = $word.ActiveDocument

# --- 1) Replace paragraph 1 content & formatting -------------------------
$p1Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t xml:space="preserve">O projeto "CX com </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cha</w:t></w:r><w:r><w:t>t</w:t></w:r><w:r><w:t>GPT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>"</w:t></w:r><w:r><w:t xml:space="preserve"> que e</w:t></w:r><w:r><w:t xml:space="preserve">stamos desenvolvendo </w:t></w:r><w:r><w:t xml:space="preserve">é </w:t></w:r><w:r><w:t xml:space="preserve">um </w:t></w:r><w:r><w:t>iniciativa</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>que visa proporcionar uma experiência de compra personalizada e satisfatória aos clientes de uma marca</w:t></w:r><w:r><w:t xml:space="preserve"> através</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>da utilização</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>de</w:t></w:r><w:r><w:t xml:space="preserve"> inteligência artificial e análise de dados para fornecer recomendações personalizadas e prever as necessidades do cliente com base em seu histórico de compras e interações anteriores com a marca. </w:t></w:r><w:r><w:t xml:space="preserve">A solução idealizada utiliza o modelo de linguagem natural </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Cha</w:t></w:r><w:r><w:t>t</w:t></w:r><w:r><w:t>GPT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, desenvolvido pela </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>OpenAI</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, para aprender com as interações anteriores dos clientes com a marca e oferecer sugestões cada vez mais precisas e personalizadas. Além disso, a solução utiliza técnicas de aprendizado de máquina para fornecer recomendações personalizadas e prever as necessidades dos clientes com base em seu histórico de compras e interações anteriores com a marca.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1Range = $d.Paragraphs(1).Range
$p1Range.InsertXML($p1Xml)

# --- 2) Insert two new paragraphs after the "Acreditamos..." paragraph ---
$acreditamosPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text
    if ($ptext -like "Acreditamos que*") {
        $acreditamosPara = $d.Paragraphs($i)
        break
    }
}

$insertRange = $acreditamosPara.Range
$insertRange.InsertParagraphAfter()
$newParaIndex = $acreditamosPara.Index + 1
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:firstLine="708"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Para o desenvolvimento da solução, serão utilizados diversos frameworks e bibliotecas Python, como o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Flask</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, para criar a aplicação web, o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TensorFlow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, para o treinamento do modelo de aprendizado de máquina, e o NLTK, para o processamento de linguagem natural.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Com relação ao modelo de receita, a solução pode ser comercializada de diversas formas, como venda de software, assinatura, comissão sobre vendas, publicidade e licenciamento de tecnologia.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs($newParaIndex).Range.InsertXML($newParaXml)
